$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.094.65"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.21%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.243.30"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.77%  "

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.68%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "117.13"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +3.07%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "265.86"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.03%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.630"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +1.22%  "

# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.41%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.608"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.66%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "46.68"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -2.30%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0927"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.51%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "9.17"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.36%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.105"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -2.42%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.41"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.08%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.883"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.00%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.585.46"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.51%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.255.43"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.97%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.339.66"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.19%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000107"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.02%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.72"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.81%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.75"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.19%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.37"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -6.33%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "231.78"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.28%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.52"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.57%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.89"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.31%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.07"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +6.01%  "

# Row 27
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.49%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "41.39"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.26%  "

# Row 29
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.72%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.24"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.28%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "172.91"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.18%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.19"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.06%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0896"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.24%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.60"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -3.24%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.35"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +10.29%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.128"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.47%  "

# Row 37
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.65"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.15%  "

# Row 38
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0372"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +4.15%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.106"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.41%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.49"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -6.98%  "

# Row 41
$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.236"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.55%  "

# Row 42
$ws.Range("B42").Value = "Celestia"
$ws.Range("C42").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "13.34"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -4.64%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "71.13"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -8.45%  "

# Row 44
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.26%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.34"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -3.62%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.63"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -9.50%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "74.68"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +32.42%  "

# Row 48
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0993"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.44%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.43"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -3.30%  "

# Row 50
$ws.Range("B50").Value = "TrustWalletToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.25"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.82%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.649"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +14.56%  "
